$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Scalpel Type" column on the (now sole) comparison row collapses from
# a set-literal of candidate strings to the single surviving type.
$ws.Range("E2").Value = "List[any]"

# Results are now committed serially, so the two extra per-comparison detail
# rows are gone and the summary rows shift up to directly follow row 2.
$ws.Rows("3:4").Delete()

# The "Total comparisons" count (now on row 3) reflects the single serial
# comparison instead of the previous batch of three.
$ws.Range("B3").Value = 1

# "Accuracy vs PyType" (now on row 5) is reported as a formatted percentage
# string rather than a bare number.
$ws.Range("F5").NumberFormat = "@"
$ws.Range("F5").Value = "100.0%"
$ws.Range("F5").NumberFormat = "general"

$wb.Save()
